$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(21.61833895679673, 21.24929000690279, 21.02601959035778, 20.93599781220204, 20.92111157509383, 21.02480146138648, 21.49047620792915, 22.42450860264021, 23.11608663681966, 23.43035241092029, 23.54918748270173, 23.5236037013353, 23.44013314406843, 23.38897926486915, 23.09553130366879, 22.9153358513931, 22.81166753930573, 22.77656679691796, 22.93452141694771, 23.4646560478722, 23.81008250476197, 23.62585676755281, 22.92584783628895, 22.17036479596043)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @(11.00307491705109, 10.48119012462744, 10.15076198695662, 10.01380894429601, 9.99093577845054, 10.14892400586205, 10.82531967089429, 12.06464558629175, 12.91339670399424, 13.28470839503007, 13.42309631377706, 13.39339203975811, 13.29613861700986, 13.23627651624429, 12.88882468983609, 12.67181494703815, 12.54560909261217, 12.50264273860912, 12.69506033947803, 13.32476521232174, 13.72334274892939, 13.51182805085805, 12.6845555882658, 11.739686020333)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colE = @(8.753575233141857, 8.720277766007317, 8.699321696015588, 8.690653368354202, 8.68920623462922, 8.699205312398629, 8.742200386213575, 8.82245870249514, 8.878958819246673, 8.904128233285231, 8.913582628154694, 8.911549868715305, 8.904907590144887, 8.900829008924038, 8.877303209280942, 8.86273436438986, 8.854304609497209, 8.851441830774796, 8.864290432907026, 8.906860671264685, 8.934235573978892, 8.91966602776475, 8.863587102308164, 8.801180124788512)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(16.86991607391233, 15.89584955866808, 15.26997757108491, 15.00819731993403, 14.96433081551593, 15.26647399323137, 16.53996406344768, 19.00274580682531, 20.67494806633232, 21.3917225636224, 21.65686569030329, 21.60004134736742, 21.4136618050453, 21.29868154950795, 20.62722412089977, 20.20408069617459, 19.95656407809808, 19.87204792380562, 20.2495528364879, 21.46857628470567, 22.22866616901555, 21.82633154475864, 20.22900810905294, 18.34778573295697)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colG = @(3.717455631414918, 3.721575550762634, 3.724234482482873, 3.725350658260761, 3.725537973554885, 3.724249403293275, 3.718849424949421, 3.709279970806527, 3.702862738942672, 3.700074800369656, 3.699037823178972, 3.699260322640221, 3.699989112476378, 3.700437955996493, 3.703047568527011, 3.704682019339776, 3.705634477537492, 3.705959091086715, 3.704506750404434, 3.699774541235019, 3.696791027564238, 3.698373429025232, 3.70458594964029, 3.711760435559823)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colG[$i]
}

$colI = @(31.72974506443047, 31.75057189724588, 31.77094388867346, 31.78114494873324, 31.78295329580638, 31.77107378466172, 31.7353475945878, 31.7257815342192, 31.75601261872514, 31.77791850831942, 31.78738953276831, 31.78529744580004, 31.77867411660802, 31.77477033147576, 31.75474539084215, 31.74455171751099, 31.7394559130079, 31.73786223413779, 31.74555740920563, 31.78058761879524, 31.81033578295263, 31.79383068360403, 31.74510035446024, 31.72185017593574)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $colI[$i]
}

$colL = @(10.36044593240418, 10.37254423197175, 10.38155955485733, 10.38563229376242, 10.38633265800336, 10.38161286623832, 10.36428800496477, 10.3429095317369, 10.33488476303186, 10.33290154261173, 10.33239008940174, 10.33248959008299, 10.33285466557767, 10.33310947382482, 10.33504790555706, 10.33666398246482, 10.33775045739746, 10.33814528041385, 10.33647570579932, 10.33274093477928, 10.33169617158703, 10.33212612391141, 10.33656033549025, 10.34734389423454)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}

$colM = @(18.02951372230606, 17.9675965194679, 17.93358872023693, 17.92074723877994, 17.91867657521113, 17.93341140725801, 18.00733805295549, 18.18368846357689, 18.33173110151102, 18.40293084264417, 18.4304314200303, 18.42448497218884, 18.40518262558179, 18.39342907132446, 18.32715428918191, 18.28747334214935, 18.26501369238949, 18.25747217721625, 18.29165990593923, 18.41083769812473, 18.49185990324279, 18.44833553355692, 18.28976605835287, 18.13268452931248)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $colM[$i]
}

$colN = @(20.68390637100179, 20.7552497753965, 20.80109839170383, 20.8202969175689, 20.82351594280936, 20.80135522392104, 20.7080820530516, 20.54134107414082, 20.42862548354468, 20.37945981057956, 20.36114440407108, 20.36507550920118, 20.37794693128791, 20.38587043075691, 20.43188095512258, 20.46064660912253, 20.47739042397224, 20.4830937196701, 20.45756391018433, 20.37415807674754, 20.32141081371033, 20.34940191659071, 20.45895695649622, 20.58472512660822)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $colN[$i]
}
